$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152 - this pushes the existing rows 152..237
# down to 153..238 (and grows the used range to A1:R238), matching the
# "weekly" price-log pattern where a new week's reading is inserted in
# chronological order rather than appended at the end.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with this week's reading for
# Zanahoria (Terminal La Palmera de La Serena / Provincia del Elquí).
$ws.Cells.Item(152, 1).Value = 8
$ws.Cells.Item(152, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(152, 3).Value = "Coquimbo"
$ws.Cells.Item(152, 4).Value = 44529
$ws.Cells.Item(152, 5).Value = 4
$ws.Cells.Item(152, 6).Value = 100114013
$ws.Cells.Item(152, 7).Value = "Zanahoria"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 600
$ws.Cells.Item(152, 11).Value = 6000
$ws.Cells.Item(152, 12).Value = 7000
$ws.Cells.Item(152, 13).Value = 6500
$ws.Cells.Item(152, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(152, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(152, 16).Value = 325
$ws.Cells.Item(152, 17).Value = 20
$ws.Cells.Item(152, 18).Value = "Hortaliza"

Write-Output "Inserted row 152; sheet now spans to row $($ws.Cells.Item(238, 4).Row)"
